$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 44797
$ws.Cells.Item(2, 10).Value = 50
$ws.Cells.Item(2, 11).Value = 10000
$ws.Cells.Item(2, 12).Value = 12000
$ws.Cells.Item(2, 13).Value = 11200
$ws.Cells.Item(2, 15).Value = 'Región de Coquimbo'
$ws.Cells.Item(2, 16).Value = 448
$ws.Cells.Item(3, 4).Value = 44482
$ws.Cells.Item(3, 10).Value = 430
$ws.Cells.Item(3, 11).Value = 8000
$ws.Cells.Item(3, 12).Value = 8500
$ws.Cells.Item(3, 13).Value = 8267
$ws.Cells.Item(3, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(3, 16).Value = 331
$ws.Cells.Item(4, 4).Value = 44497
$ws.Cells.Item(4, 10).Value = 150
$ws.Cells.Item(4, 11).Value = 6000
$ws.Cells.Item(4, 12).Value = 6500
$ws.Cells.Item(4, 13).Value = 6333
$ws.Cells.Item(4, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(4, 16).Value = 253
$ws.Cells.Item(5, 4).Value = 44799
$ws.Cells.Item(5, 10).Value = 100
$ws.Cells.Item(5, 11).Value = 10000
$ws.Cells.Item(5, 12).Value = 11000
$ws.Cells.Item(5, 13).Value = 10500
$ws.Cells.Item(5, 15).Value = 'Región de Coquimbo'
$ws.Cells.Item(5, 16).Value = 420
$ws.Cells.Item(6, 4).Value = 44537
$ws.Cells.Item(6, 10).Value = 160
$ws.Cells.Item(6, 11).Value = 8500
$ws.Cells.Item(6, 12).Value = 9000
$ws.Cells.Item(6, 13).Value = 8719
$ws.Cells.Item(6, 15).Value = 'Región del Maule'
$ws.Cells.Item(6, 16).Value = 349
$ws.Cells.Item(7, 4).Value = 44461
$ws.Cells.Item(7, 10).Value = 100
$ws.Cells.Item(7, 11).Value = 13000
$ws.Cells.Item(7, 12).Value = 14000
$ws.Cells.Item(7, 13).Value = 13500
$ws.Cells.Item(7, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(7, 16).Value = 540
$ws.Cells.Item(8, 4).Value = 44769
$ws.Cells.Item(8, 10).Value = 100
$ws.Cells.Item(8, 11).Value = 18000
$ws.Cells.Item(8, 12).Value = 20000
$ws.Cells.Item(8, 13).Value = 19000
$ws.Cells.Item(8, 15).Value = 'Región de Coquimbo'
$ws.Cells.Item(8, 16).Value = 760
$ws.Cells.Item(9, 4).Value = 44827
$ws.Cells.Item(9, 10).Value = 100
$ws.Cells.Item(9, 11).Value = 9000
$ws.Cells.Item(9, 12).Value = 10000
$ws.Cells.Item(9, 13).Value = 9500
$ws.Cells.Item(9, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(9, 16).Value = 380
$ws.Cells.Item(10, 4).Value = 44523
$ws.Cells.Item(10, 10).Value = 100
$ws.Cells.Item(10, 11).Value = 9000
$ws.Cells.Item(10, 12).Value = 10000
$ws.Cells.Item(10, 13).Value = 9500
$ws.Cells.Item(10, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(10, 16).Value = 380
$ws.Cells.Item(11, 4).Value = 44188
$ws.Cells.Item(11, 10).Value = 100
$ws.Cells.Item(11, 11).Value = 18000
$ws.Cells.Item(11, 12).Value = 20000
$ws.Cells.Item(11, 13).Value = 19000
$ws.Cells.Item(11, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(11, 16).Value = 760
$ws.Cells.Item(12, 4).Value = 44316
$ws.Cells.Item(12, 10).Value = 100
$ws.Cells.Item(12, 11).Value = 16000
$ws.Cells.Item(12, 12).Value = 18000
$ws.Cells.Item(12, 13).Value = 17000
$ws.Cells.Item(12, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(12, 16).Value = 680
$ws.Cells.Item(13, 4).Value = 44384
$ws.Cells.Item(13, 10).Value = 100
$ws.Cells.Item(13, 11).Value = 12000
$ws.Cells.Item(13, 12).Value = 13000
$ws.Cells.Item(13, 13).Value = 12500
$ws.Cells.Item(13, 15).Value = 'Región de Coquimbo'
$ws.Cells.Item(13, 16).Value = 500
$ws.Cells.Item(14, 4).Value = 44540
$ws.Cells.Item(14, 10).Value = 140
$ws.Cells.Item(14, 11).Value = 11000
$ws.Cells.Item(14, 12).Value = 12000
$ws.Cells.Item(14, 13).Value = 11429
$ws.Cells.Item(14, 15).Value = 'Región del Maule'
$ws.Cells.Item(14, 16).Value = 457
$ws.Cells.Item(15, 4).Value = 44692
$ws.Cells.Item(15, 10).Value = 100
$ws.Cells.Item(15, 11).Value = 20000
$ws.Cells.Item(15, 12).Value = 22000
$ws.Cells.Item(15, 13).Value = 21000
$ws.Cells.Item(15, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(15, 16).Value = 840
$ws.Cells.Item(16, 4).Value = 44509
$ws.Cells.Item(16, 10).Value = 100
$ws.Cells.Item(16, 11).Value = 6500
$ws.Cells.Item(16, 12).Value = 7000
$ws.Cells.Item(16, 13).Value = 6750
$ws.Cells.Item(16, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(16, 16).Value = 270
$ws.Cells.Item(17, 4).Value = 44505
$ws.Cells.Item(17, 10).Value = 180
$ws.Cells.Item(17, 11).Value = 6000
$ws.Cells.Item(17, 12).Value = 6500
$ws.Cells.Item(17, 13).Value = 6222
$ws.Cells.Item(17, 15).Value = 'Región del Maule'
$ws.Cells.Item(17, 16).Value = 249
$ws.Cells.Item(18, 4).Value = 44351
$ws.Cells.Item(18, 10).Value = 100
$ws.Cells.Item(18, 11).Value = 15000
$ws.Cells.Item(18, 12).Value = 16000
$ws.Cells.Item(18, 13).Value = 15500
$ws.Cells.Item(18, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(18, 16).Value = 620
$ws.Cells.Item(19, 4).Value = 44825
$ws.Cells.Item(19, 10).Value = 100
$ws.Cells.Item(19, 11).Value = 14000
$ws.Cells.Item(19, 12).Value = 15000
$ws.Cells.Item(19, 13).Value = 14500
$ws.Cells.Item(19, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(19, 16).Value = 580
$ws.Cells.Item(20, 4).Value = 44533
$ws.Cells.Item(20, 10).Value = 180
$ws.Cells.Item(20, 11).Value = 8000
$ws.Cells.Item(20, 12).Value = 8500
$ws.Cells.Item(20, 13).Value = 8222
$ws.Cells.Item(20, 15).Value = 'Región del Maule'
$ws.Cells.Item(20, 16).Value = 329
$ws.Cells.Item(21, 4).Value = 44813
$ws.Cells.Item(21, 10).Value = 100
$ws.Cells.Item(21, 11).Value = 11000
$ws.Cells.Item(21, 12).Value = 12000
$ws.Cells.Item(21, 13).Value = 11500
$ws.Cells.Item(21, 15).Value = 'Región de Coquimbo'
$ws.Cells.Item(21, 16).Value = 460
$ws.Cells.Item(22, 4).Value = 44503
$ws.Cells.Item(22, 10).Value = 250
$ws.Cells.Item(22, 11).Value = 9000
$ws.Cells.Item(22, 12).Value = 10000
$ws.Cells.Item(22, 13).Value = 9400
$ws.Cells.Item(22, 15).Value = 'Provincia de Melipilla'
$ws.Cells.Item(22, 16).Value = 376
$ws.Cells.Item(23, 4).Value = 44483
$ws.Cells.Item(23, 10).Value = 350
$ws.Cells.Item(23, 11).Value = 5500
$ws.Cells.Item(23, 12).Value = 6000
$ws.Cells.Item(23, 13).Value = 5714
$ws.Cells.Item(23, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(23, 16).Value = 229
$ws.Cells.Item(24, 4).Value = 44467
$ws.Cells.Item(24, 10).Value = 100
$ws.Cells.Item(24, 11).Value = 8000
$ws.Cells.Item(24, 12).Value = 9000
$ws.Cells.Item(24, 13).Value = 8500
$ws.Cells.Item(24, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(24, 16).Value = 340
$ws.Cells.Item(25, 4).Value = 44476
$ws.Cells.Item(25, 10).Value = 100
$ws.Cells.Item(25, 11).Value = 7000
$ws.Cells.Item(25, 12).Value = 7500
$ws.Cells.Item(25, 13).Value = 7250
$ws.Cells.Item(25, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(25, 16).Value = 290
$ws.Cells.Item(26, 4).Value = 44335
$ws.Cells.Item(26, 10).Value = 100
$ws.Cells.Item(26, 11).Value = 18000
$ws.Cells.Item(26, 12).Value = 20000
$ws.Cells.Item(26, 13).Value = 19000
$ws.Cells.Item(26, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(26, 16).Value = 760
$ws.Cells.Item(27, 4).Value = 44160
$ws.Cells.Item(27, 10).Value = 100
$ws.Cells.Item(27, 11).Value = 9000
$ws.Cells.Item(27, 12).Value = 10000
$ws.Cells.Item(27, 13).Value = 9500
$ws.Cells.Item(27, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(27, 16).Value = 380
$ws.Cells.Item(28, 4).Value = 44545
$ws.Cells.Item(28, 10).Value = 140
$ws.Cells.Item(28, 11).Value = 14000
$ws.Cells.Item(28, 12).Value = 15000
$ws.Cells.Item(28, 13).Value = 14429
$ws.Cells.Item(28, 15).Value = 'Provincia de Chacabuco'
$ws.Cells.Item(28, 16).Value = 577
$ws.Cells.Item(29, 4).Value = 44498
$ws.Cells.Item(29, 10).Value = 220
$ws.Cells.Item(29, 11).Value = 7000
$ws.Cells.Item(29, 12).Value = 7500
$ws.Cells.Item(29, 13).Value = 7273
$ws.Cells.Item(29, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(29, 16).Value = 291
$ws.Cells.Item(30, 4).Value = 44811
$ws.Cells.Item(30, 10).Value = 50
$ws.Cells.Item(30, 11).Value = 12000
$ws.Cells.Item(30, 12).Value = 13000
$ws.Cells.Item(30, 13).Value = 12600
$ws.Cells.Item(30, 15).Value = 'Región de Coquimbo'
$ws.Cells.Item(30, 16).Value = 504
$ws.Cells.Item(31, 4).Value = 44755
$ws.Cells.Item(31, 10).Value = 100
$ws.Cells.Item(31, 11).Value = 16000
$ws.Cells.Item(31, 12).Value = 17000
$ws.Cells.Item(31, 13).Value = 16500
$ws.Cells.Item(31, 15).Value = 'Región de Coquimbo'
$ws.Cells.Item(31, 16).Value = 660
$ws.Cells.Item(32, 4).Value = 44526
$ws.Cells.Item(32, 10).Value = 100
$ws.Cells.Item(32, 11).Value = 7500
$ws.Cells.Item(32, 12).Value = 8000
$ws.Cells.Item(32, 13).Value = 7750
$ws.Cells.Item(32, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(32, 16).Value = 310
$ws.Cells.Item(33, 4).Value = 44782
$ws.Cells.Item(33, 10).Value = 30
$ws.Cells.Item(33, 11).Value = 15000
$ws.Cells.Item(33, 12).Value = 15000
$ws.Cells.Item(33, 13).Value = 15000
$ws.Cells.Item(33, 15).Value = 'Región de Coquimbo'
$ws.Cells.Item(33, 16).Value = 600
$ws.Cells.Item(34, 4).Value = 44162
$ws.Cells.Item(34, 10).Value = 100
$ws.Cells.Item(34, 11).Value = 7500
$ws.Cells.Item(34, 12).Value = 8000
$ws.Cells.Item(34, 13).Value = 7750
$ws.Cells.Item(34, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(34, 16).Value = 310
$ws.Cells.Item(35, 4).Value = 44454
$ws.Cells.Item(35, 10).Value = 100
$ws.Cells.Item(35, 11).Value = 13000
$ws.Cells.Item(35, 12).Value = 14000
$ws.Cells.Item(35, 13).Value = 13500
$ws.Cells.Item(35, 15).Value = 'Provincia del Elquí'
$ws.Cells.Item(35, 16).Value = 540
$ws.Cells.Item(36, 4).Value = 44517
$ws.Cells.Item(36, 10).Value = 130
$ws.Cells.Item(36, 11).Value = 6000
$ws.Cells.Item(36, 12).Value = 6500
$ws.Cells.Item(36, 13).Value = 6269
$ws.Cells.Item(36, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(36, 16).Value = 251
$ws.Cells.Item(37, 4).Value = 44673
$ws.Cells.Item(37, 10).Value = 80
$ws.Cells.Item(37, 11).Value = 18000
$ws.Cells.Item(37, 12).Value = 19000
$ws.Cells.Item(37, 13).Value = 18375
$ws.Cells.Item(37, 15).Value = 'Región Metropolitana'
$ws.Cells.Item(37, 16).Value = 735
